# Weekly data update: insert a new observation row for the week of
# 2022-06-24 into the Ciboulette (Hortaliza) price series at Mercado
# Mayorista Lo Valledor de Santiago. This shifts the existing row 354
# (and everything below it) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 354, pushing old row 354..452 to 355..453.
$ws.Rows.Item(354).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(354, 1).Value = 6
$ws.Cells.Item(354, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(354, 3).Value = "Metropolitana"
$ws.Cells.Item(354, 4).Value = 44736
$ws.Cells.Item(354, 5).Value = 13
$ws.Cells.Item(354, 6).Value = 100112039
$ws.Cells.Item(354, 7).Value = "Ciboulette"
$ws.Cells.Item(354, 8).Value = "Sin especificar"
$ws.Cells.Item(354, 9).Value = "Primera"
$ws.Cells.Item(354, 10).Value = 610
$ws.Cells.Item(354, 11).Value = 700
$ws.Cells.Item(354, 12).Value = 800
$ws.Cells.Item(354, 13).Value = 741
$ws.Cells.Item(354, 14).Value = "`$/docena de atados"
$ws.Cells.Item(354, 15).Value = "Región Metropolitana"
$ws.Cells.Item(354, 16).Value = 247
$ws.Cells.Item(354, 17).Value = 3
$ws.Cells.Item(354, 18).Value = "Hortaliza"

# Match the date column's existing display style (yyyy-mm-dd hh:mm:ss).
$ws.Cells.Item(354, 4).NumberFormat = $ws.Cells.Item(355, 4).NumberFormat
